$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 55, pushing current row 55+ down to 56+.
# Excel's row insert copies the formatting of the row above into every
# column that had a styled cell there; clear the columns this new row
# doesn't actually use so only E55/F55 remain.
$ws.Rows("55").Insert()
$ws.Range("A55:D55").Clear()
$ws.Range("H55").Clear()

# Row 54 (SARIF2007 / ExpressPathsRelativeToRepoRoot): the "UriBaseId" column
# now documents the new rule name instead of "Default".
$ws.Range("E54").Value = "ProvideUriBaseIdForMappedTo"

# New row 55 documents the companion rule, still "IN PROGRESS" and styled
# the same way as the status cell directly above it (F54).
$ws.Range("E55").Value = "ExpressResultLocationsRelativeToMappedTo"
$ws.Range("F55").Value = "IN PROGRESS"
$ws.Range("F54").Copy() | Out-Null
$ws.Range("F55").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the visible view state: scrolled down, zoomed to 145%, with
# E35 selected.
$ws.Application.ActiveWindow.Zoom = 145
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("E35").Select()
